# Apply updated dSF (column F) values for specific rows.
# These reflect a repull/repush of the underlying data and a
# recalculated mean, per the commit message.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F2").Value = -6
$ws.Range("F3").Value = -6
$ws.Range("F6").Value = -5
$ws.Range("F8").Value = -8
$ws.Range("F9").Value = -3
$ws.Range("F11").Value = -3
$ws.Range("F13").Value = -8
$ws.Range("F16").Value = 1
